$d = $word.ActiveDocument

# Locate the run that currently reads "JEFA DE DEPTO. DE COMUNICACIÓN Y VINCULACIÓN."
$find = $d.Content
$found = $find.Find.Execute("JEFA DE DEPTO. DE COMUNICACIÓN Y VINCULACIÓN.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the target signature-block text to update."
}

$startPos = $find.Start

# Remove the old text entirely, then re-insert the replacement as three runs
# (JEF / E / " DE DEPTO. DE COMUNICACIÓN Y VINCULACIÓN.") via raw WordprocessingML so
# the first run keeps its original rsidRPr, matching how Word split the run when the
# "A" in "JEFA" was retyped as "E" (JEFA -> JEFE).
$find.Delete()

$target = $d.Range($startPos, $startPos)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:r w:rsidRPr="00087237">' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Montserrat Medium" w:hAnsi="Montserrat Medium"/>' +
      '<w:sz w:val="20"/>' +
      '<w:szCs w:val="20"/>' +
    '</w:rPr>' +
    '<w:t>JEF</w:t>' +
  '</w:r>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Montserrat Medium" w:hAnsi="Montserrat Medium"/>' +
      '<w:sz w:val="20"/>' +
      '<w:szCs w:val="20"/>' +
    '</w:rPr>' +
    '<w:t>E</w:t>' +
  '</w:r>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Montserrat Medium" w:hAnsi="Montserrat Medium"/>' +
      '<w:sz w:val="20"/>' +
      '<w:szCs w:val="20"/>' +
    '</w:rPr>' +
    '<w:t xml:space="preserve"> DE DEPTO. DE COMUNICACIÓN Y VINCULACIÓN.</w:t>' +
  '</w:r>' +
'</w:p>'

$target.InsertXML($xml)
